$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three stray "X cont." placeholder rows from the name list.
# (Row numbers refer to the ORIGINAL layout; delete bottom-up so earlier
# row numbers stay valid while doing so.)
$ws.Rows.Item(271).Delete()   # "D cont."
$ws.Rows.Item(137).Delete()   # "B cont."
$ws.Rows.Item(69).Delete()    # "A cont."

# Re-apply the autofilter so its range shrinks to match the new extent
# (A1:A1396) instead of the stale A1:A1399.
$ws.AutoFilterMode = $false
$ws.Range("A1:A1396").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$A`$1396"
    }
}

# Restore the scroll position / selection the sheet had when last saved.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("B64").Select()
